# Add a new "temperature + rainfall" source row (row 11) to the sources sheet,
# linking to the meteostat.net page, and adjust column C width / selection to
# match the authored change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data: id=8, data="temperature + rainfall", link -> meteostat URL
$ws.Range("B11").Value = 8
$ws.Range("C11").Value = "temperature + rainfall"

$meteoUrl = "https://meteostat.net/en/place/co/medellin?s=80110&t=2021-10-31/2024-11-16"
$ws.Hyperlinks.Add($ws.Range("D11"), $meteoUrl)
# Hyperlinks.Add() stamps extra direct font formatting on the cell that
# diverges from the existing "Hyperlink" style used by D5:D10; re-asserting
# the underline explicitly collapses it back onto the same shared style.
$ws.Range("D11").Font.Underline = 2

# Column C previously shared column B's default width; give it its own
# (wider, best-fit) width now that it holds longer text.
$ws.Range("C:C").ColumnWidth = 13.59

# Match the saved selection/active cell.
[void]$ws.Range("C11").Select()
